$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.764.00"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.407.21"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.81"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.92"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.590"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.68%  "
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.353"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.30"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.837.59"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.723.76"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.415.90"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.30"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("E19").Value = "  -1.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "328.32"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.66"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.04%  "
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("E23").Value = "  +2.79%  "
$ws.Range("E24").Value = "  +2.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.66"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.36"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0769"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.49%  "
$ws.Range("E29").Value = "  -2.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.28"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.27%  "
$ws.Range("E31").Value = "  -4.11%  "
$ws.Range("E32").Value = "  -0.48%  "
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.31"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.18"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("E38").Value = "  -2.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "320.19"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.406"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.66"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "139.65"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0966"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.57"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0513"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.95%  "
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0222"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.63%  "
$ws.Range("E48").Value = "  -3.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.53"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.05"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.56"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.32%  "
